$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-5 (columns A and B)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 109

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 85

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 77

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 44

# Remove row 6 entirely so the used range / dimension shrinks to A1:B5
$ws.Range("A6:B6").Delete()
